$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing data (A:D -> B:E)
$ws.Range("A1").EntireColumn.Insert()

# Fill the new column A with index values, styled like the header row (s=1)
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 1
$ws.Range("A4").Value = 2
$ws.Range("A5").Value = 3

$ws.Range("B1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)
